$d = $word.ActiveDocument
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    if ($r.Font.Name -ne "Times New Roman") {
        $r.Font.Name = "Times New Roman"
    }
}
